{"js": "const body = context.document.body;\n\n// 1) Replace the title text \"2.2 - Debate I\" with \"Placeholder - Check Back Later\"\nconst titleResults = body.search(\"2.2 - Debate I\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"Placeholder - Check Back Later\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Remove the trailing \" :::\" (the space run + the \":::\" run) that followed\n//    \"...general edification later.\" in the Additional Resources bullet, leaving\n//    the sentence ending in \"...general edification later.\"\nconst markerResults = body.search(\" :::\", { matchCase: true });\nmarkerResults.load(\"items\");\nawait context.sync();\n\nif (markerResults.items.length > 0) {\n  markerResults.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Replace the title text \"2.2 - Debate I\" with \"Placeholder - Check Back Later\"\n$titlePara = $d.Paragraphs.First\n$titleRange = $titlePara.Range\n$titleRange.Text = \"Placeholder - Check Back Later\"\n\n# 2) Remove the trailing \" :::\" (the space run + the \":::\" run) that followed\n#    \"...general edification later.\" in the Additional Resources bullet, without\n#    touching the preceding run's text.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \":::\"\n$found = $find.Execute()\nif ($found) {\n    $markerRange = $find.Parent\n    $deleteRange = $d.Range($markerRange.Start - 1, $markerRange.End)\n    $deleteRange.Text = \"\"\n}\n"}
